$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.686.26'
$ws.Range('E2').Value = '  -2.02%  '
$ws.Range('D3').Value = '1.810.25'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '39.35'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -9.22%  '
$ws.Range('E9').Value = '  +5.68%  '
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0991'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.04%  '
$ws.Range('D12').Value = '2.071.90'
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.674'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = '1.818.41'
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '11.15'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('E16').Value = '  -1.94%  '
$ws.Range('D17').Value = '34.690.92'
$ws.Range('E17').Value = '  -1.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').Value = '0.0₃0785'
$ws.Range('E19').Value = '  -1.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.27'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.119'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E29').Value = '  +11.14%  '
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('E31').Value = '  +2.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0545'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('E33').Value = '  -2.40%  '
$ws.Range('E34').Value = '  +16.94%  '
$ws.Range('E35').Value = '  -4.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.698'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '91.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.46%  '
$ws.Range('D39').Value = '1.323.52'
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('E40').Value = '  -0.99%  '
$ws.Range('E41').Value = '  +0.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.962'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.22'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.33%  '
$ws.Range('E45').Value = '  -4.94%  '
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0512'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.23%  '
$ws.Range('D48').Value = '1.998.06'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0668'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '98.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.20%  '
